# Applies: "Added math formula for band Fc and Q calculations"
# Target sheet: "Values" (contains the band-pass / Sallen-Key Fc & Q tables)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Values")

# --- Row 3: existing band-pass example, R3 (E3) value updated ---
$ws.Range("E3").Value = 200000

# --- Row 5: new section header "Sallen-Key low pass" (bold, like A1) ---
$ws.Range("A5").Value = "Sallen-Key low pass"
$ws.Range("A5").Font.Bold = $true

# --- Row 6: column headers for the Sallen-Key low pass table ---
$ws.Range("A6").Value = "C1 (n)"
$ws.Range("B6").Value = "C2 (n)"
$ws.Range("C6").Value = "R1"
$ws.Range("D6").Value = "R2"
$ws.Range("F6").Value = "Frequency"
$ws.Range("G6").Value = "Q"

# --- Row 7: Sallen-Key low pass example values + Fc/Q formulas ---
$ws.Range("A7").Value = 0.56
$ws.Range("B7").Value = 68
$ws.Range("C7").Value = 150000
$ws.Range("D7").Value = 150000
$ws.Range("F7").Formula = "=1/(2*PI()*SQRT(C7*`$B`$1*A7*`$B`$1*D7*B7))"
$ws.Range("G7").Formula = "=SQRT(B7/A7)/2"

# --- Row 9: new section header "Sallen-Key high pass" (bold, like A1) ---
$ws.Range("A9").Value = "Sallen-Key high pass"
$ws.Range("A9").Font.Bold = $true

# --- Row 10: column headers for the Sallen-Key high pass table ---
$ws.Range("A10").Value = "C1 (n)"
$ws.Range("B10").Value = "C2 (n)"
$ws.Range("C10").Value = "R1"
$ws.Range("D10").Value = "R2"
$ws.Range("F10").Value = "Frequency"
$ws.Range("G10").Value = "Q"

# --- Row 11: Sallen-Key high pass example values + Fc/Q formulas ---
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 3300
$ws.Range("D11").Value = 150000
$ws.Range("F11").Formula = "=1/(2*PI()*SQRT(C11*`$B`$1*A11*`$B`$1*D11*B11))"
$ws.Range("G11").Formula = "=SQRT(D11/C11)/2"

# --- Selection moved onto the newly added content ---
$ws.Range("A12").Select() | Out-Null
